$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 23666.666
$ws.Range("J93").Value = 23666.666
$ws.Range("L93").Value = 23666.666
$ws.Range("N93").Value = -28658.666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 326533.28
$ws.Range("I98").Value = 4896.8335
$ws.Range("J98").Value = 771876.0600000001
$ws.Range("K98").Value = 4896.8335
$ws.Range("L98").Value = 771876.0600000001
$ws.Range("M98").Value = -3398.8335
$ws.Range("N98").Value = -774872.0600000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 326533.28
$ws.Range("I122").Value = 4896.8335
$ws.Range("J122").Value = 771876.0600000001
$ws.Range("K122").Value = 14690.5005
$ws.Range("L122").Value = 2315628.18
$ws.Range("M122").Value = -12240.5005
$ws.Range("N122").Value = -2320528.18

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 2001.25
$ws.Range("I131").Value = 1502.5
$ws.Range("K131").Value = 4507.5
$ws.Range("M131").Value = 532.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 457443.9
$ws.Range("I32").Value = 3255.0527
$ws.Range("J32").Value = 4772238
$ws.Range("K32").Value = 3255.0527
$ws.Range("L32").Value = 4772238
$ws.Range("M32").Value = -2968.0527
$ws.Range("N32").Value = -4772812

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 28000.125
$ws.Range("J55").Value = 28000.125
$ws.Range("L55").Value = 28000.125
$ws.Range("N55").Value = -28630.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1266.6136
$ws.Range("I132").Value = 900.1
$ws.Range("J132").Value = 4931.75
$ws.Range("K132").Value = 2700.3
$ws.Range("L132").Value = 14795.25
$ws.Range("M132").Value = -170.3000000000002
$ws.Range("N132").Value = -19855.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 13701.893
$ws.Range("I82").Value = 3031.3
$ws.Range("J82").Value = 19630
$ws.Range("K82").Value = 3031.3
$ws.Range("L82").Value = 19630
$ws.Range("M82").Value = -2648.3
$ws.Range("N82").Value = -20396

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 13701.893
$ws.Range("I85").Value = 3031.3
$ws.Range("J85").Value = 19630
$ws.Range("K85").Value = 3031.3
$ws.Range("L85").Value = 19630
$ws.Range("M85").Value = -1705.3
$ws.Range("N85").Value = -22282

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3705.362
$ws.Range("I31").Value = 2353.0408
$ws.Range("J31").Value = 11068
$ws.Range("K31").Value = 2353.0408
$ws.Range("L31").Value = 11068
$ws.Range("M31").Value = -2058.0408
$ws.Range("N31").Value = -11658

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3705.362
$ws.Range("I34").Value = 2353.0408
$ws.Range("J34").Value = 11068
$ws.Range("K34").Value = 2353.0408
$ws.Range("L34").Value = 11068
$ws.Range("M34").Value = -2151.0408
$ws.Range("N34").Value = -11472

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 14428.571
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 14428.571
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 14428.571
$ws.Range("N109").Value = -16508.571
$ws.Range("M109").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 969.5
$ws.Range("I5").Value = 626
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 1878
$ws.Range("L5").Value = 6000
$ws.Range("M5").Value = -1766
$ws.Range("N5").Value = -6224

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 1125
$ws.Range("I25").Value = 250
$ws.Range("J25").Value = 2000
$ws.Range("K25").Value = 750
$ws.Range("L25").Value = 6000
$ws.Range("M25").Value = -581
$ws.Range("N25").Value = -6338

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H30").Value = 1125
$ws.Range("I30").Value = 250
$ws.Range("J30").Value = 2000
$ws.Range("K30").Value = 750
$ws.Range("L30").Value = 6000
$ws.Range("M30").Value = -648
$ws.Range("N30").Value = -6204

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 728
$ws.Range("I107").Value = 510.57144
$ws.Range("J107").Value = 897.1111
$ws.Range("K107").Value = 1531.71432
$ws.Range("L107").Value = 2691.3333
$ws.Range("M107").Value = 388.28568
$ws.Range("N107").Value = -6531.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 7463451.5
$ws.Range("I131").Value = 370.9091
$ws.Range("J131").Value = 8929414
$ws.Range("K131").Value = 1112.7273
$ws.Range("L131").Value = 26788242
$ws.Range("M131").Value = 3927.2727
$ws.Range("N131").Value = -26798322

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 969.5
$ws.Range("I135").Value = 626
$ws.Range("J135").Value = 2000
$ws.Range("K135").Value = 5634
$ws.Range("L135").Value = 18000
$ws.Range("M135").Value = -3099
$ws.Range("N135").Value = -23070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 23333.334
$ws.Range("J57").Value = 23333.334
$ws.Range("L57").Value = 23333.334
$ws.Range("N57").Value = -24973.334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 16792.58
$ws.Range("J123").Value = 16792.58
$ws.Range("L123").Value = 16792.58
$ws.Range("N123").Value = -21692.58

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1836.1936
$ws.Range("I132").Value = 1709.2543
$ws.Range("J132").Value = 4332.6665
$ws.Range("K132").Value = 5127.7629
$ws.Range("L132").Value = 12997.9995
$ws.Range("M132").Value = -2597.7629
$ws.Range("N132").Value = -18057.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2942.4614
$ws.Range("I7").Value = 2811.4443
$ws.Range("J7").Value = 3237.25
$ws.Range("K7").Value = 2811.4443
$ws.Range("L7").Value = 3237.25
$ws.Range("M7").Value = -2699.4443
$ws.Range("N7").Value = -3461.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2942.4614
$ws.Range("I126").Value = 2811.4443
$ws.Range("J126").Value = 3237.25
$ws.Range("K126").Value = 8434.332900000001
$ws.Range("L126").Value = 9711.75
$ws.Range("M126").Value = -5964.332900000001
$ws.Range("N126").Value = -14651.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 49995
$ws.Range("J109").Value = 49995
$ws.Range("L109").Value = 49995
$ws.Range("N109").Value = -52769

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 24591310
$ws.Range("I132").Value = 27273686
$ws.Range("J132").Value = 2847.1667
$ws.Range("K132").Value = 81821058
$ws.Range("L132").Value = 8541.500100000001
$ws.Range("M132").Value = -81818528
$ws.Range("N132").Value = -13601.5001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 90715
$ws.Range("J133").Value = 90715
$ws.Range("L133").Value = 90715
$ws.Range("N133").Value = -100835
